$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H41").Value = 597.6957
$ws.Range("I41").Value = 553.13336
$ws.Range("J41").Value = 681.25
$ws.Range("K41").Value = 553.13336
$ws.Range("L41").Value = 681.25
$ws.Range("M41").Value = -113.13336
$ws.Range("N41").Value = -1561.25
$ws.Range("H43").Value = 1971.0416
$ws.Range("I43").Value = 1935.7858
$ws.Range("J43").Value = 2020.4
$ws.Range("K43").Value = 1935.7858
$ws.Range("L43").Value = 2020.4
$ws.Range("M43").Value = -1866.7858
$ws.Range("N43").Value = -2158.4
$ws.Range("H53").Value = 212.7037
$ws.Range("I53").Value = 193.90909
$ws.Range("J53").Value = 225.625
$ws.Range("K53").Value = 193.90909
$ws.Range("L53").Value = 225.625
$ws.Range("M53").Value = 443.09091
$ws.Range("N53").Value = -1499.625
$ws.Range("H82").Value = 397.5
$ws.Range("I82").Value = 397.5
$ws.Range("K82").Value = 1192.5
$ws.Range("M82").Value = -786.5
$ws.Range("H85").Value = 397.5
$ws.Range("I85").Value = 397.5
$ws.Range("K85").Value = 1192.5
$ws.Range("M85").Value = 211.5
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 798
$ws.Range("I2").Value = 649.91174
$ws.Range("J2").Value = 1255.7273
$ws.Range("K2").Value = 649.91174
$ws.Range("L2").Value = 1255.7273
$ws.Range("M2").Value = -536.91174
$ws.Range("N2").Value = -1481.7273
$ws.Range("H31").Value = 1990
$ws.Range("I31").Value = 1990
$ws.Range("K31").Value = 1990
$ws.Range("M31").Value = -1696
$ws.Range("H110").Value = 1926.258
$ws.Range("I110").Value = 1378.44
$ws.Range("J110").Value = 4208.8335
$ws.Range("K110").Value = 1378.44
$ws.Range("L110").Value = 4208.8335
$ws.Range("M110").Value = 666.5599999999999
$ws.Range("N110").Value = -8298.833500000001
$ws.Range("H116").Value = 798
$ws.Range("I116").Value = 649.91174
$ws.Range("J116").Value = 1255.7273
$ws.Range("K116").Value = 649.91174
$ws.Range("L116").Value = 1255.7273
$ws.Range("M116").Value = 1644.08826
$ws.Range("N116").Value = -5843.7273
$ws.Range("H122").Value = 1593.25
$ws.Range("I122").Value = 1306.6
$ws.Range("J122").Value = 2071
$ws.Range("K122").Value = 3919.8
$ws.Range("L122").Value = 6213
$ws.Range("M122").Value = -1469.8
$ws.Range("N122").Value = -11113

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 798
$ws.Range("I3").Value = 649.91174
$ws.Range("J3").Value = 1255.7273
$ws.Range("K3").Value = 649.91174
$ws.Range("L3").Value = 1255.7273
$ws.Range("M3").Value = -535.91174
$ws.Range("N3").Value = -1483.7273
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -828
$ws.Range("N14").Value = ""
$ws.Range("H135").Value = 21874.75
$ws.Range("J135").Value = 21874.75
$ws.Range("L135").Value = 21874.75
$ws.Range("N135").Value = -32014.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 750
$ws.Range("I14").Value = 325
$ws.Range("J14").Value = 1316.6666
$ws.Range("K14").Value = 325
$ws.Range("L14").Value = 1316.6666
$ws.Range("M14").Value = -155
$ws.Range("N14").Value = -1656.6666
$ws.Range("H31").Value = 13174844
$ws.Range("I31").Value = 24391444
$ws.Range("J31").Value = 35397.086
$ws.Range("K31").Value = 24391444
$ws.Range("L31").Value = 35397.086
$ws.Range("M31").Value = -24391149
$ws.Range("N31").Value = -35987.086
$ws.Range("H34").Value = 13174844
$ws.Range("I34").Value = 24391444
$ws.Range("J34").Value = 35397.086
$ws.Range("K34").Value = 24391444
$ws.Range("L34").Value = 35397.086
$ws.Range("M34").Value = -24391242
$ws.Range("N34").Value = -35801.086
$ws.Range("H58").Value = 58824920
$ws.Range("I58").Value = 142858110
$ws.Range("J58").Value = 1682.8
$ws.Range("K58").Value = 142858110
$ws.Range("L58").Value = 1682.8
$ws.Range("M58").Value = -142857907
$ws.Range("N58").Value = -2088.8
$ws.Range("H87").Value = 29666.666
$ws.Range("J87").Value = 29666.666
$ws.Range("L87").Value = 29666.666
$ws.Range("N87").Value = -32038.666
$ws.Range("H90").Value = 29666.666
$ws.Range("J90").Value = 29666.666
$ws.Range("L90").Value = 88999.99800000001
$ws.Range("N90").Value = -100855.998
$ws.Range("H132").Value = 47626324
$ws.Range("I132").Value = 76932420
$ws.Range("J132").Value = 3929.875
$ws.Range("K132").Value = 230797260
$ws.Range("L132").Value = 11789.625
$ws.Range("M132").Value = -230794730
$ws.Range("N132").Value = -16849.625
$ws.Range("H134").Value = 4003.6
$ws.Range("I134").Value = 6216
$ws.Range("J134").Value = 2067.75
$ws.Range("K134").Value = 18648
$ws.Range("L134").Value = 6203.25
$ws.Range("M134").Value = -16113
$ws.Range("N134").Value = -11273.25
$ws.Range("H136").Value = 58824920
$ws.Range("I136").Value = 142858110
$ws.Range("J136").Value = 1682.8
$ws.Range("K136").Value = 428574330
$ws.Range("L136").Value = 5048.4
$ws.Range("M136").Value = -428571780
$ws.Range("N136").Value = -10148.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 850.7755
$ws.Range("I131").Value = 235.44444
$ws.Range("J131").Value = 989.225
$ws.Range("K131").Value = 706.33332
$ws.Range("L131").Value = 2967.675
$ws.Range("M131").Value = 4333.66668
$ws.Range("N131").Value = -13047.675
$ws.Range("H133").Value = 671143.4399999999
$ws.Range("I133").Value = 4301.8184
$ws.Range("J133").Value = 2504958
$ws.Range("K133").Value = 12905.4552
$ws.Range("L133").Value = 7514874
$ws.Range("M133").Value = -7845.4552
$ws.Range("N133").Value = -7524994
$ws.Range("H140").Value = 25811.238
$ws.Range("I140").Value = 30665.922
$ws.Range("J140").Value = 2751.5
$ws.Range("K140").Value = 91997.766
$ws.Range("L140").Value = 8254.5
$ws.Range("M140").Value = -86817.766
$ws.Range("N140").Value = -18614.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 43108.668
$ws.Range("J123").Value = 43108.668
$ws.Range("L123").Value = 43108.668
$ws.Range("N123").Value = -48008.668
$ws.Range("H126").Value = 1451.7646
$ws.Range("I126").Value = 1333.3334
$ws.Range("J126").Value = 1585
$ws.Range("K126").Value = 4000.0002
$ws.Range("L126").Value = 4755
$ws.Range("M126").Value = -1530.0002
$ws.Range("N126").Value = -9695
$ws.Range("H132").Value = 34806.09
$ws.Range("I132").Value = 44562
$ws.Range("J132").Value = 3099.375
$ws.Range("K132").Value = 133686
$ws.Range("L132").Value = 9298.125
$ws.Range("M132").Value = -131156
$ws.Range("N132").Value = -14358.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1520.8
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("H133").Value = 33007.5
$ws.Range("J133").Value = 33007.5
$ws.Range("L133").Value = 33007.5
$ws.Range("N133").Value = -43127.5
